$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("J2").Value = 1.03
$ws.Range("K2").Value = 17
$ws.Range("L2").Value = 1.17
$ws.Range("M2").Value = 5
$ws.Range("N2").Value = 1.53
$ws.Range("O2").Value = 2.4
$ws.Range("R2").Value = 1.5
$ws.Range("S2").Value = 2.5
$ws.Range("AC2").Value = 34
$ws.Range("AF2").Value = 19
# Row 4
$ws.Range("G4").Value = 5.5
$ws.Range("H4").Value = 4.1
$ws.Range("I4").Value = 1.6
$ws.Range("K4").Value = 13
$ws.Range("AF4").Value = 7.5
# Row 5
$ws.Range("G5").Value = 2.55
$ws.Range("I5").Value = 2.75
$ws.Range("U5").Value = 12
$ws.Range("Z5").Value = 9.5
$ws.Range("AE5").Value = 9
$ws.Range("AH5").Value = 29
# Row 6
$ws.Range("G6").Value = 1.75
$ws.Range("H6").Value = 3.8
$ws.Range("I6").Value = 4.5
$ws.Range("N6").Value = 2.07
$ws.Range("O6").Value = 1.83
$ws.Range("U6").Value = 8
$ws.Range("AG6").Value = 15
$ws.Range("AH6").Value = 51
$ws.Range("AI6").Value = 41
# Row 7
$ws.Range("G7").Value = 3.1
$ws.Range("I7").Value = 2.5
$ws.Range("J7").Value = 1.11
$ws.Range("K7").Value = 6.5
$ws.Range("P7").Value = 1.57
$ws.Range("Q7").Value = 2.25
$ws.Range("T7").Value = 7.5
$ws.Range("W7").Value = 34
$ws.Range("AB7").Value = 19
$ws.Range("AH7").Value = 23
$ws.Range("AI7").Value = 23
# Row 21
$ws.Range("G21").Value = 1.48
$ws.Range("I21").Value = 7.5
$ws.Range("J21").Value = 1.06
$ws.Range("K21").Value = 10
$ws.Range("R21").Value = 2.2
$ws.Range("S21").Value = 1.62
$ws.Range("T21").Value = 5.5
$ws.Range("U21").Value = 6
$ws.Range("Z21").Value = 8.5
$ws.Range("AA21").Value = 8
$ws.Range("AC21").Value = 81
# Row 23
$ws.Range("G23").Value = 3.5
$ws.Range("K23").Value = 7.5
# Row 24
$ws.Range("G24").Value = 2.63
$ws.Range("H24").Value = 2.9
$ws.Range("I24").Value = 2.9
$ws.Range("J24").Value = 1.05
$ws.Range("K24").Value = 11
$ws.Range("L24").Value = 1.2
$ws.Range("M24").Value = 4.33
$ws.Range("N24").Value = 1.75
$ws.Range("O24").Value = 2.05
$ws.Range("P24").Value = 1.33
$ws.Range("Q24").Value = 3.25
$ws.Range("R24").Value = 1.57
$ws.Range("S24").Value = 2.25
$ws.Range("T24").Value = 11
$ws.Range("U24").Value = 15
$ws.Range("V24").Value = 10
$ws.Range("W24").Value = 26
$ws.Range("Z24").Value = 11
$ws.Range("AB24").Value = 10
$ws.Range("AC24").Value = 34
$ws.Range("AD24").Value = 126
$ws.Range("AE24").Value = 12
$ws.Range("AF24").Value = 15
$ws.Range("AG24").Value = 11
$ws.Range("AH24").Value = 29
$ws.Range("AI24").Value = 21
$ws.Range("AJ24").Value = 26
# Row 25
$ws.Range("G25").Value = 3.9
$ws.Range("H25").Value = 3.2
$ws.Range("I25").Value = 2
$ws.Range("R25").Value = 1.7
$ws.Range("S25").Value = 2.05
$ws.Range("AF25").Value = 10
# Row 26
$ws.Range("L26").Value = 1.22
$ws.Range("M26").Value = 4
$ws.Range("N26").Value = 1.75
$ws.Range("O26").Value = 2.05
$ws.Range("R26").Value = 1.75
$ws.Range("S26").Value = 2
$ws.Range("U26").Value = 8.5
$ws.Range("Z26").Value = 12
$ws.Range("AA26").Value = 7.5
$ws.Range("AG26").Value = 15
# Row 27
$ws.Range("G27").Value = 3.9
$ws.Range("H27").Value = 3.9
$ws.Range("I27").Value = 1.83
$ws.Range("R27").Value = 1.7
$ws.Range("S27").Value = 2.05
$ws.Range("T27").Value = 13
$ws.Range("U27").Value = 21
$ws.Range("V27").Value = 13
$ws.Range("X27").Value = 29
$ws.Range("Y27").Value = 34
$ws.Range("AA27").Value = 7.5
$ws.Range("AB27").Value = 15
$ws.Range("AE27").Value = 8.5
$ws.Range("AF27").Value = 9.5
$ws.Range("AG27").Value = 8.5
$ws.Range("AH27").Value = 15
# Row 28
$ws.Range("G28").Value = 3.1
$ws.Range("I28").Value = 2.35
$ws.Range("V28").Value = 11
$ws.Range("X28").Value = 23
$ws.Range("Y28").Value = 29
$ws.Range("AB28").Value = 12
$ws.Range("AE28").Value = 9
$ws.Range("AF28").Value = 12
$ws.Range("AG28").Value = 9.5
$ws.Range("AH28").Value = 23
$ws.Range("AI28").Value = 19
# Row 29
$ws.Range("J29").Value = 1.03
$ws.Range("K29").Value = 15
$ws.Range("N29").Value = 1.67
$ws.Range("O29").Value = 2.15
# Row 31
$ws.Range("G31").Value = 1.95
$ws.Range("H31").Value = 3.2
$ws.Range("I31").Value = 3.8
$ws.Range("L31").Value = 1.39
$ws.Range("M31").Value = 2.57
$ws.Range("N31").Value = 2.12
$ws.Range("O31").Value = 1.57
$ws.Range("P31").Value = 1.47
$ws.Range("Q31").Value = 2.32
$ws.Range("R31").Value = 1.93
$ws.Range("S31").Value = 1.7
$ws.Range("T31").Value = 5.9
$ws.Range("U31").Value = 8.25
$ws.Range("V31").Value = 8.75
$ws.Range("W31").Value = 16.5
$ws.Range("X31").Value = 18
$ws.Range("Y31").Value = 35
$ws.Range("Z31").Value = 7.7
$ws.Range("AA31").Value = 6.2
$ws.Range("AB31").Value = 17
$ws.Range("AC31").Value = 100
$ws.Range("AD31").Value = 900
$ws.Range("AE31").Value = 9.5
$ws.Range("AF31").Value = 19.5
$ws.Range("AG31").Value = 13
$ws.Range("AH31").Value = 60
$ws.Range("AI31").Value = 40
$ws.Range("AJ31").Value = 50
# Row 32
$ws.Range("G32").Value = 2.52
$ws.Range("H32").Value = 2.9
$ws.Range("I32").Value = 2.87
$ws.Range("J32").Value = 1.13
$ws.Range("K32").Value = 4.45
$ws.Range("L32").Value = 1.6
$ws.Range("M32").Value = 2.05
$ws.Range("N32").Value = 2.72
$ws.Range("O32").Value = 1.35
$ws.Range("P32").Value = 1.62
$ws.Range("Q32").Value = 2.02
$ws.Range("R32").Value = 2.27
$ws.Range("S32").Value = 1.5
$ws.Range("T32").Value = 5.6
$ws.Range("U32").Value = 10.25
$ws.Range("V32").Value = 11.25
$ws.Range("W32").Value = 28
$ws.Range("X32").Value = 30
$ws.Range("Y32").Value = 60
$ws.Range("Z32").Value = 4.75
$ws.Range("AA32").Value = 6
$ws.Range("AB32").Value = 23
$ws.Range("AC32").Value = 175
$ws.Range("AE32").Value = 6.1
$ws.Range("AF32").Value = 12
$ws.Range("AG32").Value = 12
$ws.Range("AH32").Value = 35
$ws.Range("AI32").Value = 35
$ws.Range("AJ32").Value = 65
# Row 33
$ws.Range("T33").Value = 5.7
$ws.Range("AB33").Value = 18.5
# Row 34
$ws.Range("J34").Value = 1.05
$ws.Range("L34").Value = 1.3
$ws.Range("R34").Value = 1.92
$ws.Range("S34").Value = 1.77
# Row 35
$ws.Range("G35").Value = 4.75
$ws.Range("H35").Value = 4.2
$ws.Range("I35").Value = 1.48
$ws.Range("N35").Value = 1.22
$ws.Range("O35").Value = 4.2
$ws.Range("T35").Value = 34
$ws.Range("U35").Value = 41
$ws.Range("V35").Value = 19
$ws.Range("W35").Value = 51
$ws.Range("X35").Value = 34
$ws.Range("Y35").Value = 26
$ws.Range("AA35").Value = 12
$ws.Range("AB35").Value = 12
$ws.Range("AC35").Value = 23
$ws.Range("AD35").Value = 51
$ws.Range("AE35").Value = 19
$ws.Range("AF35").Value = 15
$ws.Range("AG35").Value = 10
$ws.Range("AH35").Value = 15
# Row 36
$ws.Range("J36").Value = 1.03
$ws.Range("L36").Value = 1.25
$ws.Range("R36").Value = 1.77
$ws.Range("S36").Value = 1.92
